$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52
$ws.Range("H52").Value = 900
$ws.Range("J52").Value = 1400
$ws.Range("L52").Value = 4200
$ws.Range("N52").Value = -4520

# Row 58
$ws.Range("H58").Value = 1170
$ws.Range("I58").Value = 320
$ws.Range("J58").Value = 2105
$ws.Range("K58").Value = 960
$ws.Range("L58").Value = 6315
$ws.Range("M58").Value = -810
$ws.Range("N58").Value = -6615

# Row 64
$ws.Range("H64").Value = 3224.7144
$ws.Range("I64").Value = 3036.2222
$ws.Range("J64").Value = 4355.6665
$ws.Range("K64").Value = 3036.2222
$ws.Range("L64").Value = 4355.6665
$ws.Range("M64").Value = -2788.2222
$ws.Range("N64").Value = -4851.6665

# Row 67
$ws.Range("H67").Value = 3224.7144
$ws.Range("I67").Value = 3036.2222
$ws.Range("J67").Value = 4355.6665
$ws.Range("K67").Value = 3036.2222
$ws.Range("L67").Value = 4355.6665
$ws.Range("M67").Value = -2178.2222
$ws.Range("N67").Value = -6071.6665

# Row 70
$ws.Range("H70").Value = 1465
$ws.Range("I70").Value = 1130
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 3390
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -3120
$ws.Range("N70").Value = -5940

# Row 73
$ws.Range("H73").Value = 1465
$ws.Range("I73").Value = 1130
$ws.Range("J73").Value = 1800
$ws.Range("K73").Value = 3390
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -2454
$ws.Range("N73").Value = -7272

# Row 74
$ws.Range("H74").Value = 3992.65
$ws.Range("J74").Value = 4019.5
$ws.Range("L74").Value = 4019.5
$ws.Range("N74").Value = -5891.5

# Row 76
$ws.Range("H76").Value = 3370462.2
$ws.Range("I76").Value = 4632941
$ws.Range("J76").Value = 3851.6667
$ws.Range("K76").Value = 4632941
$ws.Range("L76").Value = 3851.6667
$ws.Range("M76").Value = -4632626
$ws.Range("N76").Value = -4481.6667

# Row 77
$ws.Range("H77").Value = 3992.65
$ws.Range("J77").Value = 4019.5
$ws.Range("L77").Value = 20097.5
$ws.Range("N77").Value = -29457.5

# Row 79
$ws.Range("H79").Value = 3370462.2
$ws.Range("I79").Value = 4632941
$ws.Range("J79").Value = 3851.6667
$ws.Range("K79").Value = 4632941
$ws.Range("L79").Value = 3851.6667
$ws.Range("M79").Value = -4631849
$ws.Range("N79").Value = -6035.6667

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 4500
$ws.Range("I63").Value = 2928.5715
$ws.Range("J63").Value = 6333.3335
$ws.Range("K63").Value = 2928.5715
$ws.Range("L63").Value = 6333.3335
$ws.Range("M63").Value = -2242.5715
$ws.Range("N63").Value = -7705.3335

# Row 66
$ws.Range("H66").Value = 4500
$ws.Range("I66").Value = 2928.5715
$ws.Range("J66").Value = 6333.3335
$ws.Range("K66").Value = 14642.8575
$ws.Range("L66").Value = 31666.6675
$ws.Range("M66").Value = -11210.8575
$ws.Range("N66").Value = -38530.6675

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2188.0264
$ws.Range("I86").Value = 1738.8422
$ws.Range("J86").Value = 2637.2104
$ws.Range("K86").Value = 1738.8422
$ws.Range("L86").Value = 2637.2104
$ws.Range("M86").Value = -615.8422
$ws.Range("N86").Value = -4883.2104

# Row 89
$ws.Range("H89").Value = 2188.0264
$ws.Range("I89").Value = 1738.8422
$ws.Range("J89").Value = 2637.2104
$ws.Range("K89").Value = 8694.210999999999
$ws.Range("L89").Value = 13186.052
$ws.Range("M89").Value = -3078.210999999999
$ws.Range("N89").Value = -24418.052

# Row 99
$ws.Range("H99").Value = 2072.7778
$ws.Range("I99").Value = 1059.1666
$ws.Range("J99").Value = 4100
$ws.Range("K99").Value = 1059.1666
$ws.Range("L99").Value = 4100
$ws.Range("M99").Value = 438.8334
$ws.Range("N99").Value = -7096

# Row 105
$ws.Range("H105").Value = 1655.6522
$ws.Range("I105").Value = 1649
$ws.Range("K105").Value = 1649
$ws.Range("M105").Value = 98

# Row 134
$ws.Range("H134").Value = 3920.7
$ws.Range("J134").Value = 27050
$ws.Range("L134").Value = 81150
$ws.Range("N134").Value = -86220

# Row 140
$ws.Range("H140").Value = 89733.336
$ws.Range("J140").Value = 89733.336
$ws.Range("L140").Value = 89733.336
$ws.Range("N140").Value = -100093.336

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 9981.083000000001
$ws.Range("I4").Value = 9773
$ws.Range("K4").Value = 9773
$ws.Range("M4").Value = -9661

# Row 31
$ws.Range("H31").Value = 4475.7964
$ws.Range("I31").Value = 4514.9565
$ws.Range("J31").Value = 4450.778
$ws.Range("K31").Value = 4514.9565
$ws.Range("L31").Value = 4450.778
$ws.Range("M31").Value = -4219.9565
$ws.Range("N31").Value = -5040.778

# Row 34
$ws.Range("H34").Value = 4475.7964
$ws.Range("I34").Value = 4514.9565
$ws.Range("J34").Value = 4450.778
$ws.Range("K34").Value = 4514.9565
$ws.Range("L34").Value = 4450.778
$ws.Range("M34").Value = -4312.9565
$ws.Range("N34").Value = -4854.778

# Row 51
$ws.Range("H51").Value = 9350.571
$ws.Range("J51").Value = 9350.571
$ws.Range("L51").Value = 9350.571
$ws.Range("N51").Value = -10822.571

# Row 61
$ws.Range("H61").Value = 9350.571
$ws.Range("J61").Value = 9350.571
$ws.Range("L61").Value = 9350.571
$ws.Range("N61").Value = -10046.571

# Row 62
$ws.Range("H62").Value = 3305.3704
$ws.Range("I62").Value = 3333.261
$ws.Range("J62").Value = 3145
$ws.Range("K62").Value = 3333.261
$ws.Range("L62").Value = 3145
$ws.Range("M62").Value = -2709.261
$ws.Range("N62").Value = -4393

# Row 65
$ws.Range("H65").Value = 3305.3704
$ws.Range("I65").Value = 3333.261
$ws.Range("J65").Value = 3145
$ws.Range("K65").Value = 16666.305
$ws.Range("L65").Value = 15725
$ws.Range("M65").Value = -13546.305
$ws.Range("N65").Value = -21965

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 134
$ws.Range("H134").Value = 2329.5
$ws.Range("I134").Value = 1948
$ws.Range("K134").Value = 5844
$ws.Range("M134").Value = -3309

$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 22727720
$ws.Range("I18").Value = 27778152
$ws.Range("J18").Value = 780.5
$ws.Range("K18").Value = 83334456
$ws.Range("L18").Value = 2341.5
$ws.Range("M18").Value = -83334287
$ws.Range("N18").Value = -2679.5

# Row 74
$ws.Range("H74").Value = 11027.167
$ws.Range("I74").Value = 4506.5
$ws.Range("J74").Value = 14287.5
$ws.Range("K74").Value = 13519.5
$ws.Range("L74").Value = 42862.5
$ws.Range("M74").Value = -12458.5
$ws.Range("N74").Value = -44984.5

# Row 77
$ws.Range("H77").Value = 11027.167
$ws.Range("I77").Value = 4506.5
$ws.Range("J77").Value = 14287.5
$ws.Range("K77").Value = 40558.5
$ws.Range("L77").Value = 128587.5
$ws.Range("M77").Value = -35254.5
$ws.Range("N77").Value = -139195.5

# Row 81
$ws.Range("H81").Value = 2431.5
$ws.Range("J81").Value = 2431.5
$ws.Range("L81").Value = 7294.5
$ws.Range("N81").Value = -9540.5

# Row 84
$ws.Range("H84").Value = 2431.5
$ws.Range("J84").Value = 2431.5
$ws.Range("L84").Value = 21883.5
$ws.Range("N84").Value = -33115.5

# Row 107
$ws.Range("H107").Value = 398.55554
$ws.Range("I107").Value = 235.71428
$ws.Range("J107").Value = 502.18182
$ws.Range("K107").Value = 707.14284
$ws.Range("L107").Value = 1506.54546
$ws.Range("M107").Value = 1212.85716
$ws.Range("N107").Value = -5346.54546

# Row 125
$ws.Range("H125").Value = 5933.3335
$ws.Range("J125").Value = 5933.3335
$ws.Range("L125").Value = 17800.0005
$ws.Range("N125").Value = -27640.0005

# Row 126
$ws.Range("H126").Value = 3022.2222
$ws.Range("J126").Value = 3022.2222
$ws.Range("L126").Value = 9066.6666
$ws.Range("N126").Value = -18946.6666

$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 335000
$ws.Range("J44").Value = 335000
$ws.Range("L44").Value = 335000
$ws.Range("N44").Value = -336192

# Row 70
$ws.Range("H70").Value = 14430221
$ws.Range("I70").Value = 20839456
$ws.Range("J70").Value = 9443.166999999999
$ws.Range("K70").Value = 20839456
$ws.Range("L70").Value = 9443.166999999999
$ws.Range("M70").Value = -20839186
$ws.Range("N70").Value = -9983.166999999999

# Row 73
$ws.Range("H73").Value = 14430221
$ws.Range("I73").Value = 20839456
$ws.Range("J73").Value = 9443.166999999999
$ws.Range("K73").Value = 20839456
$ws.Range("L73").Value = 9443.166999999999
$ws.Range("M73").Value = -20838520
$ws.Range("N73").Value = -11315.167

# Row 80
$ws.Range("H80").Value = 2500
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1502
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 2500
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7508
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3723.9644
$ws.Range("I132").Value = 4266.1055
$ws.Range("J132").Value = 2579.4443
$ws.Range("K132").Value = 12798.3165
$ws.Range("L132").Value = 7738.3329
$ws.Range("M132").Value = -10268.3165
$ws.Range("N132").Value = -12798.3329
